# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets,
# and the single F4 change on the 演出 sheet, per the commit's regenerated
# output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2893
$ws1.Range("F15").Value = 156
$ws1.Range("F16").Value = 4158
$ws1.Range("F17").Value = 4158
$ws1.Range("F24").Value = 6218
$ws1.Range("F25").Value = 6218
$ws1.Range("F30").Value = 200
$ws1.Range("F32").Value = 5309
$ws1.Range("F36").Value = 5816
$ws1.Range("F41").Value = 3886
$ws1.Range("F42").Value = 90
$ws1.Range("F43").Value = 67
$ws1.Range("F51").Value = 2029

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 0

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 2893
$ws4.Range("F14").Value = 156
$ws4.Range("F15").Value = 4158
$ws4.Range("F16").Value = 4158
$ws4.Range("F23").Value = 6218
$ws4.Range("F24").Value = 6218
$ws4.Range("F28").Value = 200
$ws4.Range("F31").Value = 5309
$ws4.Range("F37").Value = 5816
$ws4.Range("F40").Value = 3886
$ws4.Range("F41").Value = 67
